# Apply cryptocurrency price/volume updates per upstream diff
# (GitHub Actions scheduled refresh of cryptos.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'65.112.78"
$ws.Range("E2").Value = "'  +0.99%  "

# Row 3
$ws.Range("D3").Value = "'3.379.57"
$ws.Range("E3").Value = "'  +0.59%  "

# Row 4
$ws.Range("E4").Value = "'  -0.01%  "

# Row 5
$ws.Range("D5").Value = "'554.80"
$ws.Range("E5").Value = "'  -0.21%  "

# Row 6
$ws.Range("D6").Value = "'174.16"
$ws.Range("E6").Value = "'  -0.92%  "

# Row 7
$ws.Range("D7").Value = "'0.632"
$ws.Range("E7").Value = "'  +2.01%  "

# Row 8
$ws.Range("D8").Value = "'3.369.45"
$ws.Range("E8").Value = "'  +0.47%  "

# Row 9
$ws.Range("E9").Value = "'  +0.11%  "

# Row 10
$ws.Range("D10").Value = "'0.174"
$ws.Range("E10").Value = "'  +5.59%  "

# Row 11
$ws.Range("D11").Value = "'0.637"
$ws.Range("E11").Value = "'  +1.28%  "

# Row 12
$ws.Range("D12").Value = "'53.63"
$ws.Range("E12").Value = "'  -1.76%  "

# Row 13
$ws.Range("E13").Value = "'  +1.31%  "

# Row 14
$ws.Range("D14").Value = "'9.16"
$ws.Range("E14").Value = "'  +0.66%  "

# Row 15
$ws.Range("D15").Value = "'3.914.86"
$ws.Range("E15").Value = "'  +0.82%  "

# Row 16
$ws.Range("D16").Value = "'18.33"
$ws.Range("E16").Value = "'  -0.50%  "

# Row 17
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.118"
$ws.Range("E17").Value = "'  +0.07%  "

# Row 18
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "'3.361.11"
$ws.Range("E18").Value = "'  +0.16%  "

# Row 19
$ws.Range("D19").Value = "'65.094.61"
$ws.Range("E19").Value = "'  +1.17%  "

# Row 20
$ws.Range("E20").Value = "'  -0.20%  "

# Row 21
$ws.Range("E21").Value = "'  +1.03%  "

# Row 22
$ws.Range("D22").Value = "'456.39"
$ws.Range("E22").Value = "'  -0.56%  "

# Row 23
$ws.Range("E23").Value = "'  -0.33%  "

# Row 24
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").Value = "'4.07"
$ws.Range("E24").Value = "'  -0.38%  "

# Row 25
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "'14.15"
$ws.Range("E25").Value = "'  +6.20%  "

# Row 26
$ws.Range("D26").Value = "'87.49"
$ws.Range("E26").Value = "'  +1.86%  "

# Row 27
$ws.Range("D27").Value = "'2.88"
$ws.Range("E27").Value = "'  +1.07%  "

# Row 28
$ws.Range("E28").Value = "'  -2.80%  "

# Row 29
$ws.Range("D29").Value = "'8.70"
$ws.Range("E29").Value = "'  -1.05%  "

# Row 30
$ws.Range("D30").Value = "'31.12"
$ws.Range("E30").Value = "'  +3.47%  "

# Row 31
$ws.Range("D31").Value = "'6.53"
$ws.Range("E31").Value = "'  -2.19%  "

# Row 32
$ws.Range("D32").Value = "'63.33"
$ws.Range("E32").Value = "'  +7.70%  "

# Row 33
$ws.Range("D33").Value = "'11.46"
$ws.Range("E33").Value = "'  -0.32%  "

# Row 34
$ws.Range("D34").Value = "'577.47"
$ws.Range("E34").Value = "'  -1.52%  "

# Row 35
$ws.Range("E35").Value = "'  -0.85%  "

# Row 36
$ws.Range("E36").Value = "'  -0.01%  "

# Row 37
$ws.Range("E37").Value = "'  +2.44%  "

# Row 38
$ws.Range("E38").Value = "'  +1.65%  "

# Row 39
$ws.Range("D39").Value = "'35.64"
$ws.Range("E39").Value = "'  -0.40%  "

# Row 40
$ws.Range("E40").Value = "'  -0.83%  "

# Row 41
$ws.Range("D41").Value = "'0.0₃0737"
$ws.Range("E41").Value = "'  -2.91%  "

# Row 42
$ws.Range("D42").Value = "'3.108.36"
$ws.Range("E42").Value = "'  +0.38%  "

# Row 43
$ws.Range("D43").Value = "'0.0416"
$ws.Range("E43").Value = "'  +1.06%  "

# Row 44
$ws.Range("D44").Value = "'2.75"
$ws.Range("E44").Value = "'  -1.69%  "

# Row 45
$ws.Range("D45").Value = "'3.18"
$ws.Range("E45").Value = "'  -1.07%  "

# Row 46
$ws.Range("E46").Value = "'  +2.22%  "

# Row 47
$ws.Range("E47").Value = "'  -3.70%  "

# Row 48
$ws.Range("D48").Value = "'0.999"
$ws.Range("E48").Value = "'  +0.00%  "

# Row 49
$ws.Range("D49").Value = "'140.72"
$ws.Range("E49").Value = "'  +3.77%  "

# Row 50
$ws.Range("E50").Value = "'  -2.38%  "

# Row 51
$ws.Range("D51").Value = "'8.32"
$ws.Range("E51").Value = "'  -0.59%  "

